$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new weekly record: duplicate the existing row 15 down into row 16
# (pushing the old rows 15-16 to 16-17), then overwrite row 15 with the new
# week's figures.
$ws.Rows.Item(15).Copy()
$ws.Rows.Item(16).Insert()

# New data for the newly inserted week (row 15)
$ws.Range("D15").Value = 44417
$ws.Range("J15").Value = 250
$ws.Range("K15").Value = 4000
$ws.Range("L15").Value = 4500
$ws.Range("M15").Value = 4250
$ws.Range("P15").Value = 4250
